# Applies the Valve Timing test-results update:
#  - Column E ("Errors") is narrowed from a wide text column to a slim one
#    and every per-row error message in column E is cleared.
#  - Columns C ("Open Time [s]") and D ("Close Time [s]") are populated with
#    the computed timing values for the rows that previously only had an
#    error message (now that the underlying data issue is resolved).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Narrow column E now that it no longer needs to hold long error strings.
$ws.Columns.Item(5).ColumnWidth = 5.5

$rowData = @(
    @{Row=2;  C=1.4666666666666668;  D=1.8666666666666669},
    @{Row=3;  C=2;                   D=0.80000000000000016},
    @{Row=4;  C=2.2333333333333334;  D=0.56666666666666676},
    @{Row=5;  C=2.1999999999999997;  D=1.1666666666666667},
    @{Row=6;  C=0.66666666666666663; D=2.2666666666666671},
    @{Row=7;  C=1.3;                 D=1.8333333333333333},
    @{Row=8;  C=1.3;                 D=1.7},
    @{Row=9;  C=2.4333333333333336;  D=0.66666666666666663},
    @{Row=10; C=0.80000000000000016; D=1.9666666666666668},
    @{Row=15; C=1.6333333333333335;  D=7.5333333333333341},
    @{Row=17; C=1.7666666666666668;  D=0.70000000000000007},
    @{Row=18; C=1.5;                 D=0.5},
    @{Row=21; C=1.5333333333333332;  D=0.70000000000000007},
    @{Row=24; C=1.5;                 D=0.60000000000000009},
    @{Row=26; C=1.4666666666666668;  D=6.8000000000000007},
    @{Row=29; C=0;                   D=0},
    @{Row=30; C=0.20000000000000004; D=0.20000000000000004},
    @{Row=31; C=0.23333333333333331; D=0.20000000000000004},
    @{Row=32; C=0;                   D=0},
    @{Row=33; C=0;                   D=0},
    @{Row=34; C=0;                   D=0},
    @{Row=35; C=0;                   D=0},
    @{Row=36; C=0;                   D=0},
    @{Row=37; C=0;                   D=0},
    @{Row=38; C=0;                   D=0},
    @{Row=39; C=0;                   D=0},
    @{Row=40; C=0;                   D=0},
    @{Row=41; C=0;                   D=0},
    @{Row=42; C=0;                   D=0},
    @{Row=43; C=0;                   D=0},
    @{Row=44; C=0;                   D=0}
)

foreach ($entry in $rowData) {
    $r = $entry.Row
    $ws.Cells.Item($r, 3).Value = $entry.C
    $ws.Cells.Item($r, 4).Value = $entry.D
    $ws.Cells.Item($r, 5).ClearContents()
}
